$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks on the sheet so they can be cleanly re-added
# against the refreshed F-column URLs (avoids duplicate/orphaned relationships).
$ws.Hyperlinks.Delete()

# Row 2
$ws.Cells.Item(2, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(2, 2).Value = "【急募】掲示板サイト(爆サイ)でAIによる自然な会話で書き込みを埋めていけるソフト開発者募集"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5443464"
$ws.Cells.Item(2, 7).Value = 378
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆開発 ◇サイト"

# Row 3
$ws.Cells.Item(3, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(3, 2).Value = "エッジAIカメラによる人流計測のPoC用プログラム開発"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5443336"
$ws.Cells.Item(3, 7).Value = 368
$ws.Cells.Item(3, 8).Value = "🔥AI,Ai ◆開発"

# Row 4
$ws.Cells.Item(4, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(4, 2).Value = "【急募】フロントエンド開発者募集!React/TypeScriptでのシステム構築"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5443491"
$ws.Cells.Item(4, 7).Value = 323
$ws.Cells.Item(4, 8).Value = "🔥React,TypeScript ◆開発"

# Row 5
$ws.Cells.Item(5, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(5, 2).Value = "【自動化】エクセルデータ転記作業の効率化依頼"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5442971"
$ws.Cells.Item(5, 7).Value = 145
$ws.Cells.Item(5, 8).Value = "◆効率化,自動化"

# Row 6
$ws.Cells.Item(6, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(6, 2).Value = "MT5アラートツールの制作"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "10,000 円 ~ 20,000 円 / 募集期間 3 日、取引期間 0 日"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5443470"
$ws.Cells.Item(6, 7).Value = 65
$ws.Cells.Item(6, 8).Value = "◆ツール"

# Row 7
$ws.Cells.Item(7, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(7, 2).Value = "マンション管理組合のシステム設計構築依頼"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5443592"
$ws.Cells.Item(7, 7).Value = 60
$ws.Cells.Item(7, 8).Value = "◇管理"

# Row 8
$ws.Cells.Item(8, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(8, 2).Value = "初回 2026年1月創業 コンサル会社のバックオフィス業務フロー設計・マニュアル化、IT導入 一括見積依頼"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5442904"
$ws.Cells.Item(8, 7).Value = 55
$ws.Cells.Item(8, 8).Value = "◆コンサル"

# Row 9
$ws.Cells.Item(9, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(9, 2).Value = "【急募】PGエンジニア募集!構成管理・マスタメンテ業務"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5443303"
$ws.Cells.Item(9, 7).Value = 45
$ws.Cells.Item(9, 8).Value = "◇管理"

# Row 10
$ws.Cells.Item(10, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(10, 2).Value = "【Apache Answer構築】弁護士ドットコムのような専門家Q&Aサイトのサーバー構築・初期設定"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5443617"
$ws.Cells.Item(10, 7).Value = 38
$ws.Cells.Item(10, 8).Value = "◇サイト"

# Row 11
$ws.Cells.Item(11, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(11, 2).Value = "【急募】古いPHPとPerlプログラムのアップデート依頼"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5440861"
$ws.Cells.Item(11, 7).Value = 33
$ws.Cells.Item(11, 8).Value = "○PHP"

# Row 12
$ws.Cells.Item(12, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(12, 2).Value = "【急募】Wartalesの武器アイコンとモデルを日本刀に差し替え"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5443568"
$ws.Cells.Item(12, 7).Value = 13
$ws.Cells.Item(12, 8).ClearContents()

# Row 13
$ws.Cells.Item(13, 1).Value = "2025-11-29 01:18:56"
$ws.Cells.Item(13, 2).Value = "【急募】PSE認証代行をお手伝いしてくれる方募集!"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5443188"
$ws.Cells.Item(13, 7).Value = 10
$ws.Cells.Item(13, 8).ClearContents()

# Re-create the F-column hyperlinks (row 2 through row 13) in order.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5443464")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5443336")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5443491")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5442971")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5443470")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5443592")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5442904")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5443303")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5443617")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5440861")
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5443568")
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5443188")
